$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 13.115

$ws.Range("B4").Value = 6.302999999999999
$ws.Range("D4").Value = -8.099
$ws.Range("E4").Value = 12.425

$ws.Range("D5").Value = -8.581999999999999

$ws.Range("B6").Value = 6.727999999999999
$ws.Range("D6").Value = -8.484999999999999

$ws.Range("B7").Value = 6.534000000000001

$ws.Range("B8").Value = 5.726
$ws.Range("D8").Value = -8.282999999999999

$ws.Range("E9").Value = 13.022

$ws.Range("E11").Value = 12.774

$ws.Range("E14").Value = 13.06

$ws.Range("B16").Value = 6.534000000000001
$ws.Range("D16").Value = -8.329000000000001

$ws.Range("E18").Value = 12.596

$ws.Range("B20").Value = 5.776

$ws.Range("B21").Value = 6.218999999999999

$ws.Range("D22").Value = -8.134

$ws.Range("E25").Value = 12.791
